# Insert a new data row at row 91 (pushing existing rows 91..144 down to 92..145)
# and populate it with the new record described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 91 and below down by one row.
$ws.Rows.Item(91).Insert()

# Fill in the newly inserted row with the new observation.
$ws.Range("A91").Value = 6
$ws.Range("B91").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C91").Value = "Metropolitana"
$ws.Range("D91").Value = 44572
$ws.Range("E91").Value = 13
$ws.Range("F91").Value = "Fruta"
$ws.Range("G91").Value = 100101
$ws.Range("H91").Value = "Berries"
$ws.Range("I91").Value = 100101004
$ws.Range("J91").Value = "Frambuesa"
$ws.Range("K91").Value = "Sin especificar"
$ws.Range("L91").Value = "Especial"
$ws.Range("M91").Value = 300
$ws.Range("N91").Value = 8000
$ws.Range("O91").Value = 8000
$ws.Range("P91").Value = 8000
$ws.Range("Q91").Value = "`$/bandeja 2 kilos"
$ws.Range("R91").Value = "Provincia de Linares"
$ws.Range("S91").Value = 4000
$ws.Range("T91").Value = 2
